$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.42516366666666
$ws.Range("H2").Value = 106.275491
$ws.Range("I2").Value = 0.00832770193000585
$ws.Range("J2").Value = 0.008327701930005852
$ws.Range("M2").Value = 125.901487
$ws.Range("N2").Value = 377.704461
$ws.Range("O2").Value = 0.8798726812012091
$ws.Range("P2").Value = 0.879872681201209
$ws.Range("Q2").Value = 4460.080782851705
$ws.Range("R2").Value = 40140.72704566535
$ws.Range("S2").Value = 0.007327317425398731
$ws.Range("T2").Value = 0.007327317425398732

$ws.Range("G3").Value = 35.42516366666666
$ws.Range("H3").Value = 106.275491
$ws.Range("I3").Value = 0.00832770193000585
$ws.Range("J3").Value = 0.008327701930005852
$ws.Range("O3").Value = 0.006089432091685741
$ws.Range("P3").Value = 0.006089432091685741
$ws.Range("Q3").Value = 30.86737391770122
$ws.Range("R3").Value = 277.806365259311
$ws.Range("S3").Value = 0.00005071097538257091
$ws.Range("T3").Value = 0.00005071097538257092

$ws.Range("G4").Value = 35.42516366666666
$ws.Range("H4").Value = 106.275491
$ws.Range("I4").Value = 0.00832770193000585
$ws.Range("J4").Value = 0.008327701930005852
$ws.Range("M4").Value = 2.781571666666667
$ws.Range("N4").Value = 8.344715000000001
$ws.Range("O4").Value = 0.01943923760251788
$ws.Range("P4").Value = 0.01943923760251788
$ws.Range("Q4").Value = 98.53763154222943
$ws.Range("R4").Value = 886.838683880065
$ws.Range("S4").Value = 0.0001618841765003304
$ws.Range("T4").Value = 0.0001618841765003305

$ws.Range("G5").Value = 35.42516366666666
$ws.Range("H5").Value = 106.275491
$ws.Range("I5").Value = 0.00832770193000585
$ws.Range("J5").Value = 0.008327701930005852
$ws.Range("M5").Value = 13.536175
$ws.Range("N5").Value = 40.608525
$ws.Range("O5").Value = 0.09459864910458742
$ws.Range("P5").Value = 0.09459864910458742
$ws.Range("Q5").Value = 479.5212147956416
$ws.Range("R5").Value = 4315.690933160775
$ws.Range("S5").Value = 0.0007877893527242188
$ws.Range("T5").Value = 0.000787789352724219

$ws.Range("I6").Value = 0.01070182047907406
$ws.Range("J6").Value = 0.01070182047907406
$ws.Range("M6").Value = 125.901487
$ws.Range("N6").Value = 377.704461
$ws.Range("O6").Value = 0.8798726812012091
$ws.Range("P6").Value = 0.879872681201209
$ws.Range("Q6").Value = 5731.591291502132
$ws.Range("R6").Value = 51584.32162351919
$ws.Range("S6").Value = 0.009416239478656897
$ws.Range("T6").Value = 0.009416239478656895

$ws.Range("I7").Value = 0.01070182047907406
$ws.Range("J7").Value = 0.01070182047907406
$ws.Range("O7").Value = 0.006089432091685741
$ws.Range("P7").Value = 0.006089432091685741
$ws.Range("S7").Value = 0.00006516800906473323
$ws.Range("T7").Value = 0.00006516800906473323

$ws.Range("I8").Value = 0.01070182047907406
$ws.Range("J8").Value = 0.01070182047907406
$ws.Range("M8").Value = 2.781571666666667
$ws.Range("N8").Value = 8.344715000000001
$ws.Range("O8").Value = 0.01943923760251788
$ws.Range("P8").Value = 0.01943923760251788
$ws.Range("Q8").Value = 126.6294173424317
$ws.Range("R8").Value = 1139.664756081885
$ws.Range("S8").Value = 0.0002080352310722123
$ws.Range("T8").Value = 0.0002080352310722123

$ws.Range("I9").Value = 0.01070182047907406
$ws.Range("J9").Value = 0.01070182047907406
$ws.Range("M9").Value = 13.536175
$ws.Range("N9").Value = 40.608525
$ws.Range("O9").Value = 0.09459864910458742
$ws.Range("P9").Value = 0.09459864910458742
$ws.Range("Q9").Value = 616.2264211402751
$ws.Range("R9").Value = 5546.037790262475
$ws.Range("S9").Value = 0.001012377760280214
$ws.Range("T9").Value = 0.001012377760280214

$ws.Range("G10").Value = 51.06824600000001
$ws.Range("H10").Value = 153.204738
$ws.Range("I10").Value = 0.01200505761322374
$ws.Range("J10").Value = 0.01200505761322374
$ws.Range("M10").Value = 125.901487
$ws.Range("N10").Value = 377.704461
$ws.Range("O10").Value = 0.8798726812012091
$ws.Range("P10").Value = 0.879872681201209
$ws.Range("Q10").Value = 6429.568109881804
$ws.Range("R10").Value = 57866.11298893623
$ws.Range("S10").Value = 0.01056292223012216
$ws.Range("T10").Value = 0.01056292223012216

$ws.Range("G11").Value = 51.06824600000001
$ws.Range("H11").Value = 153.204738
$ws.Range("I11").Value = 0.01200505761322374
$ws.Range("J11").Value = 0.01200505761322374
$ws.Range("O11").Value = 0.006089432091685741
$ws.Range("P11").Value = 0.006089432091685741
$ws.Range("Q11").Value = 44.49782249238868
$ws.Range("R11").Value = 400.4804024314981
$ws.Range("S11").Value = 0.00007310398309250086
$ws.Range("T11").Value = 0.00007310398309250087

$ws.Range("G12").Value = 51.06824600000001
$ws.Range("H12").Value = 153.204738
$ws.Range("I12").Value = 0.01200505761322374
$ws.Range("J12").Value = 0.01200505761322374
$ws.Range("M12").Value = 2.781571666666667
$ws.Range("N12").Value = 8.344715000000001
$ws.Range("O12").Value = 0.01943923760251788
$ws.Range("P12").Value = 0.01943923760251788
$ws.Range("Q12").Value = 142.0499861399634
$ws.Range("R12").Value = 1278.44987525967
$ws.Range("S12").Value = 0.0002333691673753724
$ws.Range("T12").Value = 0.0002333691673753725

$ws.Range("G13").Value = 51.06824600000001
$ws.Range("H13").Value = 153.204738
$ws.Range("I13").Value = 0.01200505761322374
$ws.Range("J13").Value = 0.01200505761322374
$ws.Range("M13").Value = 13.536175
$ws.Range("N13").Value = 40.608525
$ws.Range("O13").Value = 0.09459864910458742
$ws.Range("P13").Value = 0.09459864910458742
$ws.Range("Q13").Value = 691.2687147990501
$ws.Range("R13").Value = 6221.418433191451
$ws.Range("S13").Value = 0.001135662232633708
$ws.Range("T13").Value = 0.001135662232633708

$ws.Range("G14").Value = 4121.876464666667
$ws.Range("H14").Value = 12365.629394
$ws.Range("I14").Value = 0.9689654199776964
$ws.Range("J14").Value = 0.9689654199776964
$ws.Range("M14").Value = 125.901487
$ws.Range("N14").Value = 377.704461
$ws.Range("O14").Value = 0.8798726812012091
$ws.Range("P14").Value = 0.879872681201209
$ws.Range("Q14").Value = 518950.3761318364
$ws.Range("R14").Value = 4670553.385186527
$ws.Range("S14").Value = 0.8525662020670313
$ws.Range("T14").Value = 0.8525662020670312

$ws.Range("G15").Value = 4121.876464666667
$ws.Range("H15").Value = 12365.629394
$ws.Range("I15").Value = 0.9689654199776964
$ws.Range("J15").Value = 0.9689654199776964
$ws.Range("O15").Value = 0.006089432091685741
$ws.Range("P15").Value = 0.006089432091685741
$ws.Range("Q15").Value = 3591.557212681475
$ws.Range("R15").Value = 32324.01491413327
$ws.Range("S15").Value = 0.005900449124145936
$ws.Range("T15").Value = 0.005900449124145936

$ws.Range("G16").Value = 4121.876464666667
$ws.Range("H16").Value = 12365.629394
$ws.Range("I16").Value = 0.9689654199776964
$ws.Range("J16").Value = 0.9689654199776964
$ws.Range("M16").Value = 2.781571666666667
$ws.Range("N16").Value = 8.344715000000001
$ws.Range("O16").Value = 0.01943923760251788
$ws.Range("P16").Value = 0.01943923760251788
$ws.Range("Q16").Value = 11465.29478761697
$ws.Range("R16").Value = 103187.6530885527
$ws.Range("S16").Value = 0.01883594902756996
$ws.Range("T16").Value = 0.01883594902756996

$ws.Range("G17").Value = 4121.876464666667
$ws.Range("H17").Value = 12365.629394
$ws.Range("I17").Value = 0.9689654199776964
$ws.Range("J17").Value = 0.9689654199776964
$ws.Range("M17").Value = 13.536175
$ws.Range("N17").Value = 40.608525
$ws.Range("O17").Value = 0.09459864910458742
$ws.Range("P17").Value = 0.09459864910458742
$ws.Range("Q17").Value = 55794.44115410932
$ws.Range("R17").Value = 502149.9703869838
$ws.Range("S17").Value = 0.09166281975894927
$ws.Range("T17").Value = 0.09166281975894927
